$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 43; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value = 45797
}
